$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (sheet1): Property / Value table ---
$ws = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# The old sheet had two duplicate "Contact" / "No display for ContactDetail" rows
# (rows 10 and 11). Remove the duplicate row so the table collapses back down
# to a single row, then turn the remaining row into "Jurisdiction" /
# "United States of America", and fill in the previously-blank Publisher value.
$ws.Rows.Item(11).Delete()

$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# --- "Elements" sheet (sheet2): the root Extension row gets a specific
# Short/Definition instead of the generic Extension placeholder text ---
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "SSI Indicator"
$ws2.Range("L2").Value = "Indicates if the person receives Supplemental Security Income (SSI) administered via the Social Security Administration (SSA), at the time of coverage"
